$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D updates: values are text strings (e.g. "60.061.11", "4.40") that must
# remain text and not be auto-converted to numbers by Excel. Force Text number
# format before assignment, then restore the cell style to Normal afterwards so
# the cell keeps its original (unstyled) appearance.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '60.061.11'
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.422.48'
$ws.Range("D3").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '553.42'
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '137.56'
$ws.Range("D6").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '5.82'
$ws.Range("D10").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '24.98'
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.854.08'
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '60.017.77'
$ws.Range("D15").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.421.99'
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '11.42'
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.40'
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '331.95'
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.73'
$ws.Range("D21").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '65.35'
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.171'
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '8.63'
$ws.Range("D25").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.0₃0782'
$ws.Range("D28").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '170.56'
$ws.Range("D30").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '18.69'
$ws.Range("D32").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.20'
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.62'
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '39.60'
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.416'
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '313.60'
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.68'
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '139.79'
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0963'
$ws.Range("D44").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.412'
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.577'
$ws.Range("D48").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '17.77'
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '11.06'
$ws.Range("D51").Style = "Normal"

# Column E updates (percentage-like text values, e.g. "  +3.42%  ") already stay
# text since they are not parseable as plain numbers.
$ws.Range("E2").Value = '  +3.42%  '
$ws.Range("E3").Value = '  +3.11%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("E5").Value = '  +2.05%  '
$ws.Range("E6").Value = '  +2.34%  '
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("E8").Value = '  +2.33%  '
$ws.Range("E9").Value = '  +2.67%  '
$ws.Range("E10").Value = '  +5.09%  '
$ws.Range("E11").Value = '  +0.53%  '
$ws.Range("E12").Value = '  -1.97%  '
$ws.Range("E13").Value = '  +5.06%  '
$ws.Range("E14").Value = '  +3.22%  '
$ws.Range("E15").Value = '  +3.45%  '
$ws.Range("E16").Value = '  +2.77%  '
$ws.Range("E17").Value = '  +3.57%  '
$ws.Range("E18").Value = '  +6.35%  '
$ws.Range("E19").Value = '  +2.22%  '
$ws.Range("E20").Value = '  +0.87%  '
$ws.Range("E21").Value = '  +0.14%  '
$ws.Range("E22").Value = '  -0.01%  '
$ws.Range("E23").Value = '  +3.68%  '
$ws.Range("E24").Value = '  +3.91%  '
$ws.Range("E25").Value = '  +3.41%  '
$ws.Range("E26").Value = '  +0.49%  '
$ws.Range("E27").Value = '  +1.40%  '
$ws.Range("E28").Value = '  +5.94%  '
$ws.Range("E29").Value = '  +0.94%  '
$ws.Range("E30").Value = '  +0.24%  '
$ws.Range("E31").Value = '  +1.76%  '
$ws.Range("E32").Value = '  +1.88%  '
$ws.Range("E33").Value = '  +2.54%  '
$ws.Range("E35").Value = '  +5.39%  '
$ws.Range("E36").Value = '  +0.09%  '
$ws.Range("E37").Value = '  +0.63%  '
$ws.Range("E38").Value = '  +0.80%  '
$ws.Range("E39").Value = '  +1.17%  '
$ws.Range("E40").Value = '  +10.04%  '
$ws.Range("E41").Value = '  +8.17%  '
$ws.Range("E42").Value = '  +1.43%  '
$ws.Range("E43").Value = '  -0.99%  '
$ws.Range("E44").Value = '  +1.37%  '
$ws.Range("E45").Value = '  +1.74%  '
$ws.Range("E46").Value = '  +2.82%  '
$ws.Range("E47").Value = '  +7.93%  '
$ws.Range("E48").Value = '  +1.89%  '
$ws.Range("E49").Value = '  +1.67%  '
$ws.Range("E50").Value = '  +1.58%  '
$ws.Range("E51").Value = '  -0.22%  '
